$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Common Repeats" (B) / "Lone Repeats" (C) columns are stored as text
# (numeric strings), not numbers. Prefix with an apostrophe so Excel keeps
# writing them as text cells instead of converting to numeric cells.

$ws.Range("B2").Value = "'18"
$ws.Range("C2").Value = "'1"

$ws.Range("B3").Value = "'19"
$ws.Range("C3").Value = "'1"

$ws.Range("B4").Value = "'24"
$ws.Range("C4").Value = "'13"

$ws.Range("B5").Value = "'3"
$ws.Range("C5").Value = "'1"

$ws.Range("B6").Value = "'4"
$ws.Range("C6").Value = "'22"

$ws.Range("B7").Value = "'4"
$ws.Range("C7").Value = "'9"

$ws.Range("B8").Value = "'22"
$ws.Range("C8").Value = "'6"

$ws.Range("B9").Value = "'13"
$ws.Range("C9").Value = "'1"

$ws.Range("B10").Value = "'11"
$ws.Range("C10").Value = "'4"

$ws.Range("B11").Value = "'4"
$ws.Range("C11").Value = "'24"
